$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1: Update the translation version numbers:
#   3.5.1.250115      -> 3.5.1.250403CP1
#   3.5.0.230317CP4   -> 3.5.0.250403CP5
# ---------------------------------------------------------------------------

# Replace the trailing version number first (rightmost edit first so the
# forward run-coalescing performed by the engine doesn't disturb text that
# still needs to be edited).
$rngTo = $d.Content.Duplicate
$rngTo.Find.Execute("3.5.0.230317CP4", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rngTo.Text = "3.5.0.250403CP5"

$rngFrom = $d.Content.Duplicate
$rngFrom.Find.Execute("3.5.1.250115", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$rngFrom.Text = "3.5.1.250403CP1"

# The two text replacements above merge each edited span together with
# whichever neighboring runs share identical formatting, producing a single
# big run rather than the tightly-scoped runs Word would normally keep
# around a small in-place edit. Force Word to re-split run boundaries
# exactly around the two version strings by toggling (and then restoring) a
# character formatting attribute on each of them.
$probe1 = $d.Content.Duplicate
$probe1.Find.Execute("3.5.1.250403CP1", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$probe1.Bold = 1
$probe1.Bold = 0

$probe2 = $d.Content.Duplicate
$probe2.Find.Execute("3.5.0.250403CP5", $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$probe2.Bold = 1
$probe2.Bold = 0

# The bold toggle leaves behind an empty <w:rPr/> on the two freshly split
# runs even though no formatting actually changed. Clean that up by
# round-tripping the owning paragraph's OOXML through WordOpenXML /
# InsertXML with the empty <w:rPr/> markers stripped out.
$overviewPara = $d.Paragraphs(5)
$paraXml = $overviewPara.Range.WordOpenXML
$cleanXml = $paraXml -replace "<w:rPr/>", ""
$overviewPara.Range.InsertXML($cleanXml)

# ---------------------------------------------------------------------------
# Change 2: Insert a new "February 25, 2025 (...)" line (with a manual line
# break) before the existing "January 8, 2025" date line. Done after the
# version-number edits above so the Find operations there cannot accidentally
# match the version numbers that also appear in this new sentence.
# ---------------------------------------------------------------------------
$datePara = $d.Paragraphs(3)
$dateRange = $datePara.Range
$insertionPoint = $dateRange.Duplicate
$insertionPoint.Collapse(1)
$insertionPoint.InsertBefore("February 25, 2025 (Updated to versions 3.5.1.250403CP1 and 3.5.0.250403CP5; no substantive changes)")

# Re-find the start of "January" (now pushed later in the paragraph) and put
# a manual line break (<w:br/>) immediately before it, as its own run.
$searchRange = $d.Range($dateRange.Start, $dateRange.Start)
$searchRange.Find.Execute("January", $false, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$breakPoint = $d.Range($searchRange.Start, $searchRange.Start)
$breakPoint.InsertBreak(6)
